$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 43
$ws.Range("F2").Value = 29496
$ws.Range("G2").Value = 1683
$ws.Range("H2").Value = 52.58
$ws.Range("I2").Value = 44
# Row 4
$ws.Range("B4").Value = "Rocky Van Den Eeckhoudt"
$ws.Range("C4").Value = 27
$ws.Range("D4").Value = 3
$ws.Range("F4").Value = 16278
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 54.99
$ws.Range("I4").Value = 30
$ws.Range("J4").Value = 1
# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "Alessandro Delia"
$ws.Range("C5").Value = 26
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 19418
$ws.Range("G5").Value = 1138
$ws.Range("H5").Value = 51.19
$ws.Range("J5").Value = 0
# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Nick Fitzpatrick"
$ws.Range("C6").Value = 24
$ws.Range("F6").Value = 11558
$ws.Range("G6").Value = 544
$ws.Range("H6").Value = 63.74
$ws.Range("I6").Value = 26
$ws.Range("J6").Value = 1
# Row 7
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "Dartin Dan"
$ws.Range("C7").Value = 23
$ws.Range("D7").Value = 2
$ws.Range("F7").Value = 17913
$ws.Range("G7").Value = 1033
$ws.Range("H7").Value = 52.02
$ws.Range("I7").Value = 26
# Row 8
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Robin Willis"
$ws.Range("C8").Value = 19
$ws.Range("D8").Value = 0
$ws.Range("F8").Value = 14426
$ws.Range("G8").Value = 881
$ws.Range("H8").Value = 49.12
$ws.Range("I8").Value = 19
# Row 9
$ws.Range("B9").Value = "Max Walter"
$ws.Range("C9").Value = 18
$ws.Range("D9").Value = 0
$ws.Range("F9").Value = 9206
$ws.Range("G9").Value = 663
$ws.Range("H9").Value = 41.66
$ws.Range("I9").Value = 18
# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Louis Tweddle"
$ws.Range("C10").Value = 17
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 11864
$ws.Range("G10").Value = 581
$ws.Range("H10").Value = 61.26
$ws.Range("I10").Value = 18
# Row 11
$ws.Range("B11").Value = "Gijs Tromp"
$ws.Range("C11").Value = 14
$ws.Range("F11").Value = 11451
$ws.Range("G11").Value = 624
$ws.Range("H11").Value = 55.05
$ws.Range("I11").Value = 15
# Row 12
$ws.Range("B12").Value = "Milan Schoenmakers"
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = 6984
$ws.Range("G12").Value = 333
$ws.Range("H12").Value = 62.92
$ws.Range("I12").Value = 12
# Row 13
$ws.Range("A13").Value = 12
# Row 19
$ws.Range("B19").Value = "Noah B"
$ws.Range("F19").Value = 4696
$ws.Range("G19").Value = 266
$ws.Range("H19").Value = 52.96
# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "Oscar Ebbeling"
$ws.Range("C20").Value = 3
$ws.Range("F20").Value = 3730
$ws.Range("G20").Value = 258
$ws.Range("H20").Value = 43.37
$ws.Range("I20").Value = 3
# Row 21
$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Niels van Dommelen"
$ws.Range("F21").Value = 2266
$ws.Range("G21").Value = 135
$ws.Range("H21").Value = 50.36
